# Applies the "01 Neural Network Regression notes" edit:
#  - Adds <w:lang w:val="en-US"/> to a few runs and splits some runs so
#    that "Tensorflow"/"Deepmind"/"Alphafold"/"tensorflow" are wrapped in
#    <w:proofErr .../> spell-check markers (as Word's background spell
#    checker does automatically for words it doesn't recognise).
#  - Appends a new "Section 3" block with a Scikit-learn scaling link.

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    throw "No paragraph found containing: $needle"
}

function Replace-Paragraph($needle, $bodyXml) {
    $p = Get-ParagraphByText $needle
    $r = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- Paragraph 1: "Section 2 Tensorflow Fundamentals:" ---
$p1 = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
      '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Section 2</w:t></w:r>' + `
      '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:proofErr w:type="spellStart"/>' + `
      '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Tensorflow</w:t></w:r>' + `
      '<w:proofErr w:type="spellEnd"/>' + `
      '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Fundamentals</w:t></w:r>' + `
      '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r>' + `
      '</w:p>'
Replace-Paragraph "Section 2 Tensorflow Fundamentals:" $p1

# --- Paragraph: "Deepmind Alphafold" ---
$p2 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
      '<w:proofErr w:type="spellStart"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Deepmind</w:t></w:r>' + `
      '<w:proofErr w:type="spellEnd"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:proofErr w:type="spellStart"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Alphafold</w:t></w:r>' + `
      '<w:proofErr w:type="spellEnd"/>' + `
      '</w:p>'
Replace-Paragraph "Deepmind Alphafold" $p2

# --- Paragraph: "Tensorflow Tutorials:" ---
$p3 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
      '<w:proofErr w:type="spellStart"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Tensorflow</w:t></w:r>' + `
      '<w:proofErr w:type="spellEnd"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Tutorials:</w:t></w:r>' + `
      '</w:p>'
Replace-Paragraph "Tensorflow Tutorials:" $p3

# --- Paragraph: "Tensor explanation from tensorflow" ---
$p4 = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Tensor explanation from </w:t></w:r>' + `
      '<w:proofErr w:type="spellStart"/>' + `
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tensorflow</w:t></w:r>' + `
      '<w:proofErr w:type="spellEnd"/>' + `
      '</w:p>'
Replace-Paragraph "Tensor explanation from tensorflow" $p4

# --- Append the new "Section 3" block at the end of the document ---
# The document currently ends with two empty paragraphs. Collapsing to the
# very end and inserting four paragraphs means the *last* one of them
# merges into the existing trailing (empty) paragraph, while the first
# three are spliced in as genuinely new paragraphs just before it - giving
# the desired: [empty][empty][Section 3 ...][Scale ...][kaggle link].
$dash = [char]0x2013
$end = $d.Content
$end.Collapse(0)

$tailPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
    '<w:body>' + `
      '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
      '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Section 3 ' + $dash + ' Regression with </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Tensorflow</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
      '</w:p>' + `
      '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Scale, standardize or norma</w:t></w:r>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>lize with scikit learn</w:t></w:r>' + `
      '</w:p>' + `
      '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:hyperlink r:id="rIdScikit" w:history="1">' + `
          '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:lang w:val="en-US"/></w:rPr>' + `
          '<w:t>https://www.kaggle.com/discdiver/guide-to-scaling-and-standardizing</w:t></w:r>' + `
        '</w:hyperlink>' + `
      '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part>' + `
    '<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">' + `
    '<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' + `
      '<Relationship Id="rIdScikit" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.kaggle.com/discdiver/guide-to-scaling-and-standardizing" TargetMode="External"/>' + `
    '</Relationships></pkg:xmlData></pkg:part>' + `
    '</pkg:package>'

$end.InsertXML($tailPkg)

Write-Output "done"
